# Insert a new daily-ranking row for 2026/01/23 (金) right before the
# existing row 707, shifting the 2026/12/29 .. 2027/01/05 block down by one
# row (707 -> 708, ..., 748 -> 749) and growing the used range from
# A1:D748 to A1:D749.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 707 (and everything below it) down by one row, leaving a blank
# row 707 behind for the new record.
$ws.Rows.Item(707).Insert()

# Column A holds plain text dates (e.g. "2026/12/29"), stored as literal
# strings rather than date serials. Writing a date-shaped string directly
# into ".Value" gets auto-recognized and converted to a real date by
# Excel's smart typing, so format the cell as Text first, then clear the
# formatting again afterwards (ClearFormats keeps the already-stored
# string value intact while dropping the now-unneeded explicit style, so
# the cell matches the unstyled look of every other data row).
$ws.Range("A707").NumberFormat = "@"
$ws.Range("A707").Value = "2026/01/23"
$ws.Range("A707").ClearFormats()

$ws.Range("B707").Value = "金"
$ws.Range("C707").Value = 3
$ws.Range("D707").Value = 174
